# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 775
$ws1.Range("F6").Value = 130
$ws1.Range("F8").Value = 142
$ws1.Range("F9").Value = 332
$ws1.Range("F10").Value = 444
$ws1.Range("F11").Value = 504
$ws1.Range("F13").Value = 11600
$ws1.Range("F14").Value = 5395

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 775
$ws4.Range("F8").Value = 130
$ws4.Range("F10").Value = 142
$ws4.Range("F11").Value = 332
$ws4.Range("F12").Value = 444
$ws4.Range("F13").Value = 504
$ws4.Range("F15").Value = 11600
$ws4.Range("F17").Value = 5395
